# Apply "Error Calculations and Plots" edits to missing_data.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill/clear a few imputed values in the upper block (rows 19, 21, 23) ---
$ws.Range("C19").Value = 13.2
$ws.Range("C21").ClearContents()
$ws.Range("C23").Value = 12.2

# --- Remove two rows that no longer belong in the data set ---
# Delete from the bottom up so row indices of the earlier row stay valid.
$ws.Rows.Item(28).Delete()   # old row 28 "SC 92"
$ws.Rows.Item(26).Delete()   # old row 26 "RM 232"

# --- After the deletions, the two former "SC 101" / "SC 232" rows (now at
#     rows 27 and 33) also have their Column C values swapped ---
$ws.Range("C27").ClearContents()
$ws.Range("C33").Value = 10.4
